$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "South Creek" conductivity row (row 23), mirroring the
# existing Conductivity row (21) but with a new "Conductivity" label.
$ws.Range("A23").Value = "Conductivity"
$ws.Range("B23").Value = "COND"
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "uS/cm"

# Match the new selection left in the sheet after the edit.
$ws.Range("H1:H20").Select()
